# Update Data by bot, scripted by HH
#
# Applies the row-2 data refresh described by the commit diff:
#   J2, M2, N2            -> text fields (report/date codes)
#   O2,P2,Q2,S2,U2,W2,X2,
#   Z2,AB2,AF2,AG2         -> numeric balance-sheet figures
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference cell (K2) keeps its original, unmodified default style - use it
# to restore style on the text cells below after forcing a text number
# format, so we don't leave a stray explicit style on the edited cells.
$defaultStyle = $ws.Range("K2").Style

# --- Text (string) fields -------------------------------------------------
# Force text format before assigning so values like "001" keep their
# leading zero and the date-like strings are stored as plain text instead
# of being auto-converted to numbers/dates.
$ws.Range("J2").NumberFormat = "@"
$ws.Range("J2").Value = "001"
$ws.Range("J2").Style = $defaultStyle

$ws.Range("M2").NumberFormat = "@"
$ws.Range("M2").Value = "2020-12-17 00:00:00"
$ws.Range("M2").Style = $defaultStyle

$ws.Range("N2").NumberFormat = "@"
$ws.Range("N2").Value = "2017-12-31 00:00:00"
$ws.Range("N2").Style = $defaultStyle

# --- Numeric fields --------------------------------------------------------
$ws.Range("O2").Value = 2028641985.53
$ws.Range("P2").Value = 296691140.95
$ws.Range("Q2").Value = 53032608.98
$ws.Range("S2").Value = 234409028.67
$ws.Range("U2").Value = 325966282.69
$ws.Range("W2").Value = 1727574863.67
$ws.Range("X2").Value = 134245249
$ws.Range("Z2").Value = 125910479.04
$ws.Range("AB2").Value = 301067121.86
$ws.Range("AF2").Value = 106.8957754648
$ws.Range("AG2").Value = 85.1591791944
